$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 1.53  # G2: 1.5 -> 1.53
$ws.Cells.Item(2, 9).Value = 6.25  # I2: 6.5 -> 6.25
$ws.Cells.Item(2, 12).Value = 6  # L2: 6.5 -> 6
$ws.Cells.Item(2, 17).Value = 1.85  # Q2: 1.88 -> 1.85
$ws.Cells.Item(2, 18).Value = 2  # R2: 1.98 -> 2
$ws.Cells.Item(2, 19).Value = 1.36  # S2: 1.4 -> 1.36
$ws.Cells.Item(2, 20).Value = 3  # T2: 2.75 -> 3
$ws.Cells.Item(2, 21).Value = 1.91  # U2: 2 -> 1.91
$ws.Cells.Item(2, 22).Value = 1.8  # V2: 1.73 -> 1.8
$ws.Cells.Item(2, 24).Value = 7  # X2: 6.5 -> 7
$ws.Cells.Item(2, 26).Value = 11  # Z2: 10 -> 11
$ws.Cells.Item(2, 32).Value = 51  # AF2: 67 -> 51
$ws.Cells.Item(2, 35).Value = 29  # AI2: 34 -> 29
$ws.Cells.Item(2, 36).Value = 19  # AJ2: 21 -> 19
$ws.Cells.Item(2, 38).Value = 41  # AL2: 51 -> 41
$ws.Cells.Item(2, 46).Value = 3  # AT2: 2.75 -> 3
$ws.Cells.Item(2, 48).Value = 51  # AV2: 67 -> 51
$ws.Cells.Item(2, 54).Value = 301  # BB2: 351 -> 301
$ws.Cells.Item(3, 7).Value = 2.45  # G3: 2.4 -> 2.45
$ws.Cells.Item(3, 9).Value = 3.2  # I3: 3.25 -> 3.2
$ws.Cells.Item(3, 12).Value = 4  # L3: 4.33 -> 4
$ws.Cells.Item(3, 19).Value = 1.67  # S3: 1.62 -> 1.67
$ws.Cells.Item(3, 20).Value = 2.1  # T3: 2.2 -> 2.1
$ws.Cells.Item(3, 35).Value = 13  # AI3: 15 -> 13
$ws.Cells.Item(3, 46).Value = 2.1  # AT3: 2.2 -> 2.1
$ws.Cells.Item(3, 47).Value = 10  # AU3: 9.5 -> 10
$ws.Cells.Item(3, 52).Value = 67  # AZ3: 81 -> 67
$ws.Cells.Item(4, 7).Value = 3.5  # G4: 3.2 -> 3.5
$ws.Cells.Item(4, 8).Value = 3  # H4: 2.88 -> 3
$ws.Cells.Item(4, 9).Value = 2.3  # I4: 2.45 -> 2.3
$ws.Cells.Item(4, 10).Value = 4  # J4: 3.75 -> 4
$ws.Cells.Item(4, 11).Value = 1.95  # K4: 2 -> 1.95
$ws.Cells.Item(4, 12).Value = 3.1  # L4: 3.25 -> 3.1
$ws.Cells.Item(4, 15).Value = 1.44  # O4: 1.4 -> 1.44
$ws.Cells.Item(4, 16).Value = 2.63  # P4: 2.75 -> 2.63
$ws.Cells.Item(4, 17).Value = 2.4  # Q4: 2.3 -> 2.4
$ws.Cells.Item(4, 18).Value = 1.53  # R4: 1.6 -> 1.53
$ws.Cells.Item(4, 19).Value = 1.53  # S4: 1.5 -> 1.53
$ws.Cells.Item(4, 20).Value = 2.38  # T4: 2.5 -> 2.38
$ws.Cells.Item(4, 21).Value = 2  # U4: 1.91 -> 2
$ws.Cells.Item(4, 22).Value = 1.73  # V4: 1.8 -> 1.73
$ws.Cells.Item(4, 25).Value = 13  # Y4: 12 -> 13
$ws.Cells.Item(4, 26).Value = 41  # Z4: 34 -> 41
$ws.Cells.Item(4, 27).Value = 34  # AA4: 29 -> 34
$ws.Cells.Item(4, 30).Value = 6  # AD4: 5.5 -> 6
$ws.Cells.Item(4, 31).Value = 17  # AE4: 15 -> 17
$ws.Cells.Item(4, 32).Value = 67  # AF4: 51 -> 67
$ws.Cells.Item(4, 34).Value = 6.5  # AH4: 7 -> 6.5
$ws.Cells.Item(4, 35).Value = 9.5  # AI4: 11 -> 9.5
$ws.Cells.Item(4, 37).Value = 21  # AK4: 23 -> 21
$ws.Cells.Item(4, 40).Value = 5  # AN4: 4.75 -> 5
$ws.Cells.Item(4, 41).Value = 21  # AO4: 17 -> 21
$ws.Cells.Item(4, 42).Value = 34  # AP4: 29 -> 34
$ws.Cells.Item(4, 43).Value = 67  # AQ4: 51 -> 67
$ws.Cells.Item(4, 44).Value = 101  # AR4: 81 -> 101
$ws.Cells.Item(4, 45).Value = 301  # AS4: 251 -> 301
$ws.Cells.Item(4, 46).Value = 2.38  # AT4: 2.5 -> 2.38
$ws.Cells.Item(4, 47).Value = 9  # AU4: 8.5 -> 9
$ws.Cells.Item(4, 49).Value = 4  # AW4: 4.33 -> 4
$ws.Cells.Item(4, 50).Value = 13  # AX4: 15 -> 13
$ws.Cells.Item(4, 51).Value = 29  # AY4: 26 -> 29
$ws.Cells.Item(4, 52).Value = 41  # AZ4: 51 -> 41
$ws.Cells.Item(5, 7).Value = 3.25  # G5: 3.1 -> 3.25
$ws.Cells.Item(5, 9).Value = 2.25  # I5: 2.35 -> 2.25
$ws.Cells.Item(5, 10).Value = 4  # J5: 3.75 -> 4
$ws.Cells.Item(5, 17).Value = 2.25  # Q5: 2.3 -> 2.25
$ws.Cells.Item(5, 18).Value = 1.62  # R5: 1.6 -> 1.62
$ws.Cells.Item(5, 23).Value = 8.5  # W5: 8 -> 8.5
$ws.Cells.Item(5, 33).Value = 1250  # AG5: 1000 -> 1250
$ws.Cells.Item(5, 34).Value = 6.5  # AH5: 7 -> 6.5
$ws.Cells.Item(5, 36).Value = 9.5  # AJ5: 10 -> 9.5
$ws.Cells.Item(5, 52).Value = 41  # AZ5: 51 -> 41
$ws.Cells.Item(5, 53).Value = 67  # BA5: 81 -> 67
